$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91; this shifts the existing rows 91..129
# down to 92..130, carrying all their original values (including styles)
# with them automatically.
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new record's data. Columns
# not explicitly called out below (A, B, C, E, F, G, H, I, N, O, Q, R) keep
# the same values the prior occupant of row 91 had, since Insert() shifted
# that whole row down intact and we only overwrite the changed fields here.
$ws.Cells.Item(91, 1).Value = 4
$ws.Cells.Item(91, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(91, 3).Value = "Los Lagos"
$ws.Cells.Item(91, 4).Value = 44274
$ws.Cells.Item(91, 4).NumberFormat = $ws.Cells.Item(92, 4).NumberFormat
$ws.Cells.Item(91, 5).Value = 10
$ws.Cells.Item(91, 6).Value = 100112032
$ws.Cells.Item(91, 7).Value = "Zapallo italiano"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 150
$ws.Cells.Item(91, 11).Value = 1200
$ws.Cells.Item(91, 12).Value = 1200
$ws.Cells.Item(91, 13).Value = 1200
$ws.Cells.Item(91, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(91, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(91, 16).Value = 20
$ws.Cells.Item(91, 17).Value = 60
$ws.Cells.Item(91, 18).Value = "Hortaliza"
